$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.71305630770718
$ws.Range("C2").Value = 9.451336989984247
$ws.Range("D2").Value = 13.50359275942293
$ws.Range("E2").Value = 13.93488704166072
$ws.Range("G2").Value = 34.7162323040215
$ws.Range("H2").Value = 15.87222554054175
$ws.Range("I2").Value = 24.35881136032112
$ws.Range("J2").Value = 8.465716674283758
$ws.Range("L2").Value = 12.39292078574622
$ws.Range("O2").Value = 24.93698879047697

$ws.Range("B3").Value = 18.09721428897859
$ws.Range("C3").Value = 9.136980873044592
$ws.Range("D3").Value = 13.51811972828297
$ws.Range("E3").Value = 13.98145022051006
$ws.Range("G3").Value = 34.87907222383711
$ws.Range("H3").Value = 15.93962130545373
$ws.Range("I3").Value = 24.51451632222741
$ws.Range("J3").Value = 8.4779571114834
$ws.Range("L3").Value = 12.36379497826929
$ws.Range("O3").Value = 25.05474597235135

$ws.Range("B4").Value = 17.70872994337077
$ws.Range("C4").Value = 8.937647344400389
$ws.Range("D4").Value = 13.52962895136157
$ws.Range("E4").Value = 14.01203978604673
$ws.Range("G4").Value = 34.99132140972459
$ws.Range("H4").Value = 15.98397169866864
$ws.Range("I4").Value = 24.61559136546436
$ws.Range("J4").Value = 8.485882445163252
$ws.Range("L4").Value = 12.34737732359793
$ws.Range("O4").Value = 25.13317750405693

$ws.Range("B5").Value = 17.54802995317773
$ws.Range("C5").Value = 8.854923076194245
$ws.Range("D5").Value = 13.53497023772901
$ws.Range("E5").Value = 14.02500855782237
$ws.Range("G5").Value = 35.04013079196577
$ws.Range("H5").Value = 16.00279121490691
$ws.Range("I5").Value = 24.65815721629309
$ws.Range("J5").Value = 8.489215400279839
$ws.Range("L5").Value = 12.34105976767248
$ws.Range("O5").Value = 25.16667637289257

$ws.Range("B6").Value = 17.52120816782281
$ws.Range("C6").Value = 8.841099466840028
$ws.Range("D6").Value = 13.53589648066944
$ws.Range("E6").Value = 14.02719242284715
$ws.Range("G6").Value = 35.04842027607533
$ws.Range("H6").Value = 16.00596125361138
$ws.Range("I6").Value = 24.66530844088981
$ws.Range("J6").Value = 8.489775084055347
$ws.Range("L6").Value = 12.34003337473962
$ws.Range("O6").Value = 25.17233158370829

$ws.Range("B7").Value = 17.70657206454101
$ws.Range("C7").Value = 8.936537614568142
$ws.Range("D7").Value = 13.52969834932341
$ws.Range("E7").Value = 14.01221264907276
$ws.Range("G7").Value = 34.99196727533839
$ws.Range("H7").Value = 15.98422248353167
$ws.Range("I7").Value = 24.61615984673097
$ws.Range("J7").Value = 8.485926975827459
$ws.Range("L7").Value = 12.34729060812188
$ws.Range("O7").Value = 25.13362306147916

$ws.Range("B8").Value = 18.5029857563935
$ws.Range("C8").Value = 9.344311982391512
$ws.Range("D8").Value = 13.50806424370579
$ws.Range("E8").Value = 13.95052729073031
$ws.Range("G8").Value = 34.7698237842658
$ws.Range("H8").Value = 15.89484721423458
$ws.Range("I8").Value = 24.41136344899028
$ws.Range("J8").Value = 8.469852337453689
$ws.Range("L8").Value = 12.3825765954379
$ws.Range("O8").Value = 24.9763174175981

$ws.Range("B9").Value = 19.97414764963374
$ws.Range("C9").Value = 10.09014638297749
$ws.Range("D9").Value = 13.48618042059781
$ws.Range("E9").Value = 13.84540736717522
$ws.Range("G9").Value = 34.43225955643792
$ws.Range("H9").Value = 15.74315186932185
$ws.Range("I9").Value = 24.05312340819814
$ws.Range("J9").Value = 8.441566499650442
$ws.Range("L9").Value = 12.4632079352849
$ws.Range("O9").Value = 24.71662539855321

$ws.Range("B10").Value = 20.98979249411325
$ws.Range("C10").Value = 10.60099087597231
$ws.Range("D10").Value = 13.48260562714255
$ws.Range("E10").Value = 13.7778062258152
$ws.Range("G10").Value = 34.24500695040529
$ws.Range("H10").Value = 15.64608433544694
$ws.Range("I10").Value = 23.81629212780031
$ws.Range("J10").Value = 8.422738417668032
$ws.Range("L10").Value = 12.52914588598653
$ws.Range("O10").Value = 24.55578721481716

$ws.Range("B11").Value = 21.43592006900214
$ws.Range("C11").Value = 10.82457735567136
$ws.Range("D11").Value = 13.48368609086871
$ws.Range("E11").Value = 13.74913788290039
$ws.Range("G11").Value = 34.1732086574138
$ws.Range("H11").Value = 15.60505253838623
$ws.Range("I11").Value = 23.71426395061338
$ws.Range("J11").Value = 8.414593077402909
$ws.Range("L11").Value = 12.56053410630199
$ws.Range("O11").Value = 24.48916932945477

$ws.Range("B12").Value = 21.60244701901156
$ws.Range("C12").Value = 10.90792563095063
$ws.Range("D12").Value = 13.48448334782627
$ws.Range("E12").Value = 13.73858115237892
$ws.Range("G12").Value = 34.14795997673773
$ws.Range("H12").Value = 15.5899646105011
$ws.Range("I12").Value = 23.67644864430818
$ws.Range("J12").Value = 8.411568680393666
$ws.Range("L12").Value = 12.57261453582962
$ws.Range("O12").Value = 24.46488840440643

$ws.Range("B13").Value = 21.56669158182704
$ws.Range("C13").Value = 10.8900345401178
$ws.Range("D13").Value = 13.48429440305074
$ws.Range("E13").Value = 13.74084142413365
$ws.Range("G13").Value = 34.15331122721877
$ws.Range("H13").Value = 15.59319404559721
$ws.Range("I13").Value = 23.68455635483711
$ws.Range("J13").Value = 8.412217371333535
$ws.Range("L13").Value = 12.57000424943929
$ws.Range("O13").Value = 24.47007560378563

$ws.Range("B14").Value = 21.44966934921386
$ws.Range("C14").Value = 10.83146118888783
$ws.Range("D14").Value = 13.48374391053738
$ws.Range("E14").Value = 13.74826337657696
$ws.Range("G14").Value = 34.17109247334925
$ws.Range("H14").Value = 15.60380222487413
$ws.Range("I14").Value = 23.71113641827347
$ws.Range("J14").Value = 8.414343056203194
$ws.Range("L14").Value = 12.56152410904371
$ws.Range("O14").Value = 24.48715274349748

$ws.Range("B15").Value = 21.37767226310731
$ws.Range("C15").Value = 10.79541005633202
$ws.Range("D15").Value = 13.48345722452481
$ws.Range("E15").Value = 13.75284851226098
$ws.Range("G15").Value = 34.18223705207835
$ws.Range("H15").Value = 15.61035865223809
$ws.Range("I15").Value = 23.72752433061088
$ws.Range("J15").Value = 8.41565291378339
$ws.Range("L15").Value = 12.55635491559343
$ws.Range("O15").Value = 24.49773628279869

$ws.Range("B16").Value = 20.96030603062687
$ws.Range("C16").Value = 10.58619716179853
$ws.Range("D16").Value = 13.48258939185803
$ws.Range("E16").Value = 13.7797216747067
$ws.Range("G16").Value = 34.24996977334596
$ws.Range("H16").Value = 15.64882878730523
$ws.Range("I16").Value = 23.82307474118238
$ws.Range("J16").Value = 8.423279155361509
$ws.Range("L16").Value = 12.52712211934693
$ws.Range("O16").Value = 24.56027297641623

$ws.Range("B17").Value = 20.70009972624019
$ws.Range("C17").Value = 10.45555772810021
$ws.Range("D17").Value = 13.48274952989629
$ws.Range("E17").Value = 13.79674097808581
$ws.Range("G17").Value = 34.29496060966462
$ws.Range("H17").Value = 15.67322969811502
$ws.Range("I17").Value = 23.88315349788384
$ws.Range("J17").Value = 8.428064890793923
$ws.Range("L17").Value = 12.50954133979614
$ws.Range("O17").Value = 24.60031736044599

$ws.Range("B18").Value = 20.54894643285378
$ws.Range("C18").Value = 10.3795919190613
$ws.Range("D18").Value = 13.48309648074625
$ws.Range("E18").Value = 13.80672618110142
$ws.Range("G18").Value = 34.32209704413334
$ws.Range("H18").Value = 15.68755853232831
$ws.Range("I18").Value = 23.91824643707002
$ws.Range("J18").Value = 8.430857035042226
$ws.Range("L18").Value = 12.4995607953068
$ws.Range("O18").Value = 24.62396592988464

$ws.Range("B19").Value = 20.49751679794964
$ws.Range("C19").Value = 10.35373117615679
$ws.Range("D19").Value = 13.48325775193381
$ws.Range("E19").Value = 13.81014069726635
$ws.Range("G19").Value = 34.33150071426775
$ws.Range("H19").Value = 15.69246051945643
$ws.Range("I19").Value = 23.9302205801547
$ws.Range("J19").Value = 8.431809202510564
$ws.Range("L19").Value = 12.49620430751721
$ws.Range("O19").Value = 24.63207861647926

$ws.Range("B20").Value = 20.72795425291269
$ws.Range("C20").Value = 10.46955036169542
$ws.Range("D20").Value = 13.48270611484977
$ws.Range("E20").Value = 13.79490894556126
$ws.Range("G20").Value = 34.29004085458027
$ws.Range("H20").Value = 15.6706017395274
$ws.Range("I20").Value = 23.87670241234085
$ws.Range("J20").Value = 8.427551353373236
$ws.Range("L20").Value = 12.51139928060583
$ws.Range("O20").Value = 24.59599077297704

$ws.Range("B21").Value = 21.48410800553763
$ws.Range("C21").Value = 10.8487017966665
$ws.Range("D21").Value = 13.48389508030564
$ws.Range("E21").Value = 13.74607524728374
$ws.Range("G21").Value = 34.1658169286328
$ws.Range("H21").Value = 15.6006741298533
$ws.Range("I21").Value = 23.70330694082908
$ws.Range("J21").Value = 8.413717063046885
$ws.Range("L21").Value = 12.5640097040162
$ws.Range("O21").Value = 24.48211107093419

$ws.Range("B22").Value = 21.96419165342771
$ws.Range("C22").Value = 11.08879025058265
$ws.Range("D22").Value = 13.48693365951168
$ws.Range("E22").Value = 13.71590430651731
$ws.Range("G22").Value = 34.09594232616791
$ws.Range("H22").Value = 15.55759518140655
$ws.Range("I22").Value = 23.59476575780671
$ws.Range("J22").Value = 8.405025541062022
$ws.Range("L22").Value = 12.5995238457428
$ws.Range("O22").Value = 24.41319900734852

$ws.Range("B23").Value = 21.70929063017337
$ws.Range("C23").Value = 10.96137195899647
$ws.Range("D23").Value = 13.48510540711866
$ws.Range("E23").Value = 13.73184755931717
$ws.Range("G23").Value = 34.13219579189961
$ws.Range("H23").Value = 15.58034704020196
$ws.Range("I23").Value = 23.65225864045094
$ws.Range("J23").Value = 8.409632437490796
$ws.Range("L23").Value = 12.58046786433404
$ws.Range("O23").Value = 24.44947267848486

$ws.Range("B24").Value = 20.71536606206584
$ws.Range("C24").Value = 10.46322696495568
$ws.Range("D24").Value = 13.48272494877727
$ws.Range("E24").Value = 13.79573658265826
$ws.Range("G24").Value = 34.292261119141
$ws.Range("H24").Value = 15.67178890382326
$ws.Range("I24").Value = 23.87961722593325
$ws.Range("J24").Value = 8.427783396650568
$ws.Range("L24").Value = 12.510558909909
$ws.Range("O24").Value = 24.59794487189214

$ws.Range("B25").Value = 19.58696046818442
$ws.Range("C25").Value = 9.894617150965455
$ws.Range("D25").Value = 13.48990266810024
$ws.Range("E25").Value = 13.87215182100711
$ws.Range("G25").Value = 34.51297836111939
$ws.Range("H25").Value = 15.78166480533849
$ws.Range("I25").Value = 24.14540177277642
$ws.Range("J25").Value = 8.448874107748797
$ws.Range("L25").Value = 12.4401979948655
$ws.Range("O25").Value = 24.78163382156595
